$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.749.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.35%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.350.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.09%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.83%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.01%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.349.69'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.17%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.475'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.35%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.55'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.90%  '

# Row 11
$ws.Range('E11').Value = '  +2.04%  '

# Row 12
$ws.Range('E12').Value = '  +2.50%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.923.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.89%  '

# Row 14
$ws.Range('E14').Value = '  +1.46%  '

# Row 15
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000171'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.27%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.349.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.89%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.83%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.850.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.27%  '

# Row 19
$ws.Range('E19').Value = '  +5.46%  '

# Row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.17%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '371.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '

# Row 23
$ws.Range('E23').Value = '  +2.67%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.487.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.01%  '

# Row 25
$ws.Range('E25').Value = '  -0.02%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.73%  '

# Row 27
$ws.Range('E27').Value = '  +9.77%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.94%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.83%  '

# Row 30
$ws.Range('E30').Value = '  +0.37%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.30%  '

# Row 32
$ws.Range('E32').Value = '  +0.43%  '

# Row 33
$ws.Range('E33').Value = '  +2.44%  '

# Row 34
$ws.Range('E34').Value = '  -0.07%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.383.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.03%  '

# Row 36
$ws.Range('E36').Value = '  +2.21%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.50'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.80%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.93%  '

# Row 39
$ws.Range('E39').Value = '  +3.40%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.42'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.22%  '

# Row 41
$ws.Range('E41').Value = '  +3.29%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.20%  '

# Row 43
$ws.Range('E43').Value = '  -1.43%  '

# Row 44
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.20%  '

# Row 45
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.754'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.19%  '

# Row 46
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.14%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.36%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.97%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.16%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.86%  '

# Row 51
$ws.Range('E51').Value = '  +12.62%  '
